$d = $word.ActiveDocument
$replacements = @(
    @("71×82=5822", "31×63=1953"),
    @("42×14=588", "54×95=5130"),
    @("45×14=630", "99×32=3168"),
    @("39×11=429", "13×42=546"),
    @("53×86=4558", "51×32=1632"),
    @("87×51=4437", "73×93=6789"),
    @("47×56=2632", "25×42=1050"),
    @("32×26=832", "69×55=3795"),
    @("45×55=2475", "84×58=4872"),
    @("78×28=2184", "47×72=3384"),
    @("99×61=6039", "29×12=348"),
    @("65×85=5525", "60×24=1440"),
    @("37×70=2590", "43×33=1419"),
    @("17×26=442", "55×52=2860"),
    @("13×19=247", "12×17=204"),
    @("58×61=3538", "31×54=1674"),
    @("49×21=1029", "36×70=2520"),
    @("84×86=7224", "15×88=1320"),
    @("92×74=6808", "47×49=2303"),
    @("41×79=3239", "77×65=5005"),
    @("64×36=2304", "11×61=671"),
    @("35×45=1575", "92×90=8280"),
    @("79×67=5293", "59×89=5251"),
    @("81×65=5265", "26×27=702"),
    @("42×40=1680", "70×14=980"),
)
foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}